$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database values for rows 2 and 3 (columns D:AQ)
$ws.Range("D2").Value = 0.166
$ws.Range("D3").Value = 0.166

$ws.Range("E2").Value = 0.113
$ws.Range("E3").Value = 0.113

$ws.Range("F2").Value = 0.04190000000000001
$ws.Range("F3").Value = 0.04190000000000001

$ws.Range("G2").Value = 0.05740956520600673
$ws.Range("G3").Value = 0.05740956520600673

$ws.Range("H2").Value = 0.05740956520600673
$ws.Range("H3").Value = 0.05740956520600673

$ws.Range("I2").Value = 0.04347939932705778
$ws.Range("I3").Value = 0.04347939932705778

$ws.Range("J2").Value = 0.04347939932705778
$ws.Range("J3").Value = 0.04347939932705778

$ws.Range("K2").Value = 230.2
$ws.Range("K3").Value = 230.2

$ws.Range("L2").Value = 0.03013838520050798
$ws.Range("L3").Value = 0.03013838520050798

$ws.Range("M2").Value = 8.43
$ws.Range("M3").Value = 8.43

$ws.Range("N2").Value = 0.002407952240852352
$ws.Range("N3").Value = 0.002407952240852352

$ws.Range("O2").Value = 0.03662033014769765
$ws.Range("O3").Value = 0.03662033014769765

$ws.Range("P2").Value = 1.17
$ws.Range("P3").Value = 1.17

$ws.Range("Q2").Value = 0.0003341997772001485
$ws.Range("Q3").Value = 0.0003341997772001485

$ws.Range("R2").Value = 0.005082536924413553
$ws.Range("R3").Value = 0.005082536924413553

$ws.Range("S2").Value = 7.26
$ws.Range("S3").Value = 7.26

$ws.Range("T2").Value = 0.8612099644128114
$ws.Range("T3").Value = 0.8612099644128114

$ws.Range("U2").Value = 412.7
$ws.Range("U3").Value = 412.7

$ws.Range("V2").Value = 0.1178839726927361
$ws.Range("V3").Value = 0.1178839726927361

$ws.Range("W2").Value = 0.06407972386148535
$ws.Range("W3").Value = 0.06407972386148535

$ws.Range("X2").Value = 0.09048629046646217
$ws.Range("X3").Value = 0.09048629046646217

$ws.Range("Y2").Value = -0.02640656660497682
$ws.Range("Y3").Value = -0.02640656660497682

$ws.Range("Z2").Value = 1.211281677186083
$ws.Range("Z3").Value = 1.211281677186083

$ws.Range("AA2").Value = 0.05266579973992198
$ws.Range("AA3").Value = 0.05266579973992198

$ws.Range("AB2").Value = 0.05546958694069165
$ws.Range("AB3").Value = 0.05546958694069165

$ws.Range("AC2").Value = -0.002803787200769674
$ws.Range("AC3").Value = -0.002803787200769674

$ws.Range("AD2").Value = 3514.4
$ws.Range("AD3").Value = 3514.4

$ws.Range("AE2").Value = 0
$ws.Range("AE3").Value = 0

$ws.Range("AF2").Value = 3514.4
$ws.Range("AF3").Value = 3514.4

$ws.Range("AG2").Value = 3101.7
$ws.Range("AG3").Value = 3101.7

$ws.Range("AH2").Value = 0.5009621826579048
$ws.Range("AH3").Value = 0.5009621826579048

$ws.Range("AI2").Value = 0.4833514420498151
$ws.Range("AI3").Value = 0.4833514420498151

$ws.Range("AJ2").Value = 0.4697694847484324
$ws.Range("AJ3").Value = 0.4697694847484324

$ws.Range("AK2").Value = 0.4522615263480214
$ws.Range("AK3").Value = 0.4522615263480214

$ws.Range("AL2").Value = 94.8
$ws.Range("AL3").Value = 94.8

$ws.Range("AM2").Value = 94.8
$ws.Range("AM3").Value = 94.8

$ws.Range("AN2").Value = 10.05838580423583
$ws.Range("AN3").Value = 10.05838580423583

$ws.Range("AO2").Value = 3.503164556962026
$ws.Range("AO3").Value = 3.503164556962026

$ws.Range("AP2").Value = 8.877218088151118
$ws.Range("AP3").Value = 8.877218088151118

$ws.Range("AQ2").Value = 3.503164556962026
$ws.Range("AQ3").Value = 3.503164556962026

